$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2964644314608904
$ws.Range("C2").Value = 0.05204127957692606
$ws.Range("D2").Value = 0.3057664949759413
$ws.Range("F2").Value = 1.286094136379496
$ws.Range("G2").Value = 0.6339111078817652
$ws.Range("H2").Value = 0.7629503473669246
$ws.Range("J2").Value = 0.3173351099104806
$ws.Range("K2").Value = 0.2717952620865276
$ws.Range("M2").Value = 0.2516613191299228
$ws.Range("O2").Value = 2.781774390993235
$ws.Range("B3").Value = 0.2625494012835929
$ws.Range("C3").Value = 0.04888545475904493
$ws.Range("D3").Value = 0.2999429382849428
$ws.Range("F3").Value = 1.291720718292183
$ws.Range("G3").Value = 0.6387555901064914
$ws.Range("H3").Value = 0.7689265224063533
$ws.Range("J3").Value = 0.3160395519624686
$ws.Range("K3").Value = 0.2373991465323257
$ws.Range("M3").Value = 0.2375135800254284
$ws.Range("O3").Value = 2.804307867632033
$ws.Range("B4").Value = 0.2417079151751977
$ws.Range("C4").Value = 0.04694026051144817
$ws.Range("D4").Value = 0.2964919325168012
$ws.Range("F4").Value = 1.295844628252006
$ws.Range("G4").Value = 0.6421300859656895
$ws.Range("H4").Value = 0.7729054917629057
$ws.Range("J4").Value = 0.31540878542706
$ws.Range("K4").Value = 0.2162145717714168
$ws.Range("M4").Value = 0.2289072308901936
$ws.Range("O4").Value = 2.819630914032317
$ws.Range("B5").Value = 0.2332109866188148
$ws.Range("C5").Value = 0.04614575421661016
$ws.Range("D5").Value = 0.295117115070255
$ws.Range("F5").Value = 1.29769352697334
$ws.Range("G5").Value = 0.6436057574119971
$ws.Range("H5").Value = 0.7746048606151845
$ws.Range("J5").Value = 0.3151931961716485
$ws.Range("K5").Value = 0.2075658426622766
$ws.Range("M5").Value = 0.2254205263047027
$ws.Range("O5").Value = 2.826249243200451
$ws.Range("B6").Value = 0.231799861621937
$ws.Range("C6").Value = 0.04601371876749027
$ws.Range("D6").Value = 0.2948907341396421
$ws.Range("F6").Value = 1.298010707525059
$ws.Range("G6").Value = 0.6438568626223216
$ws.Range("H6").Value = 0.7748917466975982
$ws.Range("J6").Value = 0.3151599024507519
$ws.Range("K6").Value = 0.2061287873534354
$ws.Range("M6").Value = 0.224842803733722
$ws.Range("O6").Value = 2.827370805531345
$ws.Range("B7").Value = 0.241593337496397
$ws.Range("C7").Value = 0.04692955282185096
$ws.Range("D7").Value = 0.2964732635332012
$ws.Range("F7").Value = 1.295868881276689
$ws.Range("G7").Value = 0.6421495803663575
$ws.Range("H7").Value = 0.7729280945231878
$ws.Range("J7").Value = 0.3154057100256082
$ws.Range("K7").Value = 0.2160979953928859
$ws.Range("M7").Value = 0.2288601248553235
$ws.Range("O7").Value = 2.8197186565083
$ws.Range("B8").Value = 0.284774485755463
$ws.Range("C8").Value = 0.05095474580920722
$ws.Range("D8").Value = 0.3037327384858202
$ws.Range("F8").Value = 1.28789535762867
$ws.Range("G8").Value = 0.6354984499891643
$ws.Range("H8").Value = 0.76494671761359
$ws.Range("J8").Value = 0.3168542356635484
$ws.Range("K8").Value = 0.2599493415768848
$ws.Range("M8").Value = 0.2467666209462251
$ws.Range("O8").Value = 2.789235249605071
$ws.Range("B9").Value = 0.3692934290116909
$ws.Range("C9").Value = 0.05878618197132823
$ws.Range("D9").Value = 0.3189525309677208
$ws.Range("F9").Value = 1.27756511855285
$ws.Range("G9").Value = 0.6256311423882366
$ws.Range("H9").Value = 0.7517488909760459
$ws.Range("J9").Value = 0.3210006965727459
$ws.Range("K9").Value = 0.3454037864930513
$ws.Range("M9").Value = 0.2825110815476179
$ws.Range("O9").Value = 2.741258097520131
$ws.Range("B10").Value = 0.4312714502480901
$ws.Range("C10").Value = 0.06449949078040618
$ws.Range("D10").Value = 0.3307286920296093
$ws.Range("F10").Value = 1.273206379858159
$ws.Range("G10").Value = 0.62032076248083
$ws.Range("H10").Value = 0.7435446351446586
$ws.Range("J10").Value = 0.3248426742733272
$ws.Range("K10").Value = 0.4078378036124946
$ws.Range("M10").Value = 0.309148280072101
$ws.Range("O10").Value = 2.713202991964053
$ws.Range("B11").Value = 0.459437145465273
$ws.Range("C11").Value = 0.06708931097215043
$ws.Range("D11").Value = 0.3362138741861429
$ws.Range("F11").Value = 1.271924343341425
$ws.Range("G11").Value = 0.618326722494146
$ws.Range("H11").Value = 0.7401355696988219
$ws.Range("J11").Value = 0.3267631752027285
$ws.Range("K11").Value = 0.4361606297782714
$ws.Range("M11").Value = 0.3213462402573697
$ws.Range("O11").Value = 2.70200208721576
$ws.Range("B12").Value = 0.4700981965796132
$ws.Range("C12").Value = 0.06806863053625989
$ws.Range("D12").Value = 0.3383092715922942
$ws.Range("F12").Value = 1.271539561413007
$ws.Range("G12").Value = 0.6176323209658676
$ws.Range("H12").Value = 0.7388910487844527
$ws.Range("J12").Value = 0.3275152400289016
$ws.Range("K12").Value = 0.4468739435208988
$ws.Range("M12").Value = 0.3259766792292993
$ws.Range("O12").Value = 2.69798513721372
$ws.Range("B13").Value = 0.4678023664379793
$ws.Range("C13").Value = 0.06785777939373361
$ws.Range("D13").Value = 0.3378571797347547
$ws.Range("F13").Value = 1.271617953714006
$ws.Range("G13").Value = 0.6177791722266193
$ws.Range("H13").Value = 0.7391570148357687
$ws.Range("J13").Value = 0.3273521663499395
$ws.Range("K13").Value = 0.4445671797505213
$ws.Range("M13").Value = 0.3249789311326268
$ws.Range("O13").Value = 2.698840269852042
$ws.Range("B14").Value = 0.4603143339112705
$ws.Range("C14").Value = 0.06716990840590142
$ws.Range("D14").Value = 0.3363858983302919
$ws.Range("F14").Value = 1.271890669395688
$ws.Range("G14").Value = 0.6182683768333561
$ws.Range("H14").Value = 0.7400322521645109
$ws.Range("J14").Value = 0.3268245509206622
$ws.Range("K14").Value = 0.4370422630223061
$ws.Range("M14").Value = 0.3217269634906827
$ws.Range("O14").Value = 2.701667108947134
$ws.Range("B15").Value = 0.4557270696869296
$ws.Range("C15").Value = 0.06674838466966548
$ws.Range("D15").Value = 0.3354870718800953
$ws.Range("F15").Value = 1.272070827221349
$ws.Range("G15").Value = 0.6185759352741087
$ws.Range("H15").Value = 0.740574404074863
$ws.Range("J15").Value = 0.326504601676632
$ws.Range("K15").Value = 0.4324314632058588
$ws.Range("M15").Value = 0.3197365083269403
$ws.Range("O15").Value = 2.703427879081261
$ws.Range("B16").Value = 0.4294301326719392
$ws.Range("C16").Value = 0.06433004920175733
$ws.Range("D16").Value = 0.3303727883569962
$ws.Range("F16").Value = 1.273304259119051
$ws.Range("G16").Value = 0.6204595681907819
$ws.Range("H16").Value = 0.7437739234236034
$ws.Range("J16").Value = 0.3247206388953856
$ws.Range("K16").Value = 0.4059852044071306
$ws.Range("M16").Value = 0.3083527145789517
$ws.Range("O16").Value = 2.713966422942633
$ws.Range("B17").Value = 0.4132901078120312
$ws.Range("C17").Value = 0.0628440762884992
$ws.Range("D17").Value = 0.3272680564035539
$ws.Range("F17").Value = 1.274240363276192
$ws.Range("G17").Value = 0.6217231618229988
$ws.Range("H17").Value = 0.7458194464818746
$ws.Range("J17").Value = 0.3236704632755476
$ws.Range("K17").Value = 0.3897406871305691
$ws.Range("M17").Value = 0.3013895833259213
$ws.Range("O17").Value = 2.720831450803743
$ws.Range("B18").Value = 0.4040041533827718
$ws.Range("C18").Value = 0.06198852270966881
$ws.Range("D18").Value = 0.3254943674930075
$ws.Range("F18").Value = 1.274844750373106
$ws.Range("G18").Value = 0.6224896342208339
$ws.Range("H18").Value = 0.7470263906560888
$ws.Range("J18").Value = 0.3230826934507007
$ws.Range("K18").Value = 0.3803898969030968
$ws.Range("M18").Value = 0.2973921724392028
$ws.Range("O18").Value = 2.724927006611068
$ws.Range("B19").Value = 0.400859653819623
$ws.Range("C19").Value = 0.06169870110139186
$ws.Range("D19").Value = 0.3248959046248814
$ws.Range("F19").Value = 1.275060717185042
$ws.Range("G19").Value = 0.6227559629641846
$ws.Range("H19").Value = 0.7474402660332231
$ws.Range("J19").Value = 0.3228864788751338
$ws.Range("K19").Value = 0.3772226320789969
$ws.Range("M19").Value = 0.2960400303542414
$ws.Range("O19").Value = 2.726338932793766
$ws.Range("B20").Value = 0.41500851849068
$ws.Range("C20").Value = 0.0630023502676238
$ws.Range("D20").Value = 0.3275973122214992
$ws.Range("F20").Value = 1.274133886814028
$ws.Range("G20").Value = 0.6215845422457136
$ws.Range("H20").Value = 0.7455985496906621
$ws.Range("J20").Value = 0.3237805732825478
$ws.Range("K20").Value = 0.3914707112549536
$ws.Range("M20").Value = 0.3021300359696681
$ws.Range("O20").Value = 2.72008544532946
$ws.Range("B21").Value = 0.4625138825385022
$ws.Range("C21").Value = 0.06737199102366276
$ws.Range("D21").Value = 0.336817554476653
$ws.Range("F21").Value = 1.271807833876544
$ws.Range("G21").Value = 0.6181230377675604
$ws.Range("H21").Value = 0.7397739143021482
$ws.Range("J21").Value = 0.3269788511960172
$ws.Range("K21").Value = 0.4392528419540724
$ws.Range("M21").Value = 0.3226818395396549
$ws.Range("O21").Value = 2.700830702116491
$ws.Range("B22").Value = 0.4935337908244151
$ws.Range("C22").Value = 0.07021967599706613
$ws.Range("D22").Value = 0.3429499842716837
$ws.Range("F22").Value = 1.270874531168161
$ws.Range("G22").Value = 0.616214567529525
$ws.Range("H22").Value = 0.7362377197268017
$ws.Range("J22").Value = 0.3292137241887332
$ws.Range("K22").Value = 0.4704114338692875
$ws.Range("M22").Value = 0.3361795915948562
$ws.Range("O22").Value = 2.689555714973181
$ws.Range("B23").Value = 0.476980604646684
$ws.Range("C23").Value = 0.06870057874488111
$ws.Range("D23").Value = 0.3396672986446561
$ws.Range("F23").Value = 1.271318974744972
$ws.Range("G23").Value = 0.6172007591519417
$ws.Range("H23").Value = 0.7381003125924082
$ws.Range("J23").Value = 0.3280077085179727
$ws.Range("K23").Value = 0.4537880926557989
$ws.Range("M23").Value = 0.3289696341895194
$ws.Range("O23").Value = 2.695453589764497
$ws.Range("B24").Value = 0.4142316462805695
$ws.Range("C24").Value = 0.06293079848408922
$ws.Range("D24").Value = 0.3274484205734325
$ws.Range("F24").Value = 1.274181818583074
$ws.Range("G24").Value = 0.6216470875155196
$ws.Range("H24").Value = 0.7456983207322452
$ws.Range("J24").Value = 0.3237307427389879
$ws.Range("K24").Value = 0.3906886034384627
$ws.Range("M24").Value = 0.3017952592886317
$ws.Range("O24").Value = 2.720422251030499
$ws.Range("B25").Value = 0.3464480333491906
$ws.Range("C25").Value = 0.05667448216657078
$ws.Range("D25").Value = 0.3147304185840341
$ws.Range("F25").Value = 1.279792075684917
$ws.Range("G25").Value = 0.627960210070249
$ws.Range("H25").Value = 0.7550569492478161
$ws.Range("J25").Value = 0.3197392177416631
$ws.Range("K25").Value = 0.3223459063559915
$ws.Range("M25").Value = 0.2727746779665807
$ws.Range("O25").Value = 2.752973870889477
